$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

Set-TextValue "D2" "46.206.56"
Set-TextValue "E2" "  -1.18%  "

Set-TextValue "D3" "2.503.45"
Set-TextValue "E3" "  +10.78%  "

Set-TextValue "D4" "0.996"
Set-TextValue "E4" "  -0.41%  "

Set-TextValue "D5" "297.29"
Set-TextValue "E5" "  -0.17%  "

Set-TextValue "D6" "97.36"
Set-TextValue "E6" "  -0.97%  "

Set-TextValue "D7" "0.581"
Set-TextValue "E7" "  +1.94%  "

Set-TextValue "D8" "0.999"
Set-TextValue "E8" "  -0.13%  "

Set-TextValue "D9" "0.534"
Set-TextValue "E9" "  +6.44%  "

Set-TextValue "D10" "36.11"
Set-TextValue "E10" "  +4.72%  "

Set-TextValue "D11" "0.0795"
Set-TextValue "E11" "  +2.06%  "

Set-TextValue "D12" "7.44"
Set-TextValue "E12" "  +6.76%  "

Set-TextValue "B13" "TRON"
Set-TextValue "C13" "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue "D13" "0.104"
Set-TextValue "E13" "  +2.25%  "

Set-TextValue "B14" "WrappedliquidstakedEther2.0"
Set-TextValue "C14" "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue "D14" "2.862.77"
Set-TextValue "E14" "  +9.96%  "

Set-TextValue "D15" "2.493.63"
Set-TextValue "E15" "  +10.51%  "

Set-TextValue "D16" "0.873"
Set-TextValue "E16" "  +10.48%  "

Set-TextValue "D17" "14.47"
Set-TextValue "E17" "  +7.49%  "

Set-TextValue "D18" "46.176.93"
Set-TextValue "E18" "  -1.17%  "

Set-TextValue "D19" "13.10"
Set-TextValue "E19" "  +6.29%  "

Set-TextValue "D20" "0.0₃0959"
Set-TextValue "E20" "  -0.56%  "

Set-TextValue "D21" "6.42"
Set-TextValue "E21" "  +11.32%  "

Set-TextValue "D22" "68.19"
Set-TextValue "E22" "  +4.05%  "

Set-TextValue "D23" "249.13"
Set-TextValue "E23" "  +2.16%  "

Set-TextValue "D24" "2.84"
Set-TextValue "E24" "  +3.02%  "

Set-TextValue "D25" "2.00"
Set-TextValue "E25" "  +8.53%  "

Set-TextValue "D26" "0.999"
Set-TextValue "E26" "  -0.16%  "

Set-TextValue "D27" "40.71"
Set-TextValue "E27" "  +0.35%  "

Set-TextValue "D28" "2.23"
Set-TextValue "E28" "  +0.73%  "

Set-TextValue "D29" "9.97"
Set-TextValue "E29" "  +5.02%  "

Set-TextValue "D30" "22.19"
Set-TextValue "E30" "  +11.49%  "

Set-TextValue "D31" "3.95"
Set-TextValue "E31" "  +19.58%  "

Set-TextValue "B32" "Filecoin"
Set-TextValue "C32" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D32" "5.70"
Set-TextValue "E32" "  +7.84%  "

Set-TextValue "B33" "ARBITRUM"
Set-TextValue "C33" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D33" "2.19"
Set-TextValue "E33" "  +33.17%  "

Set-TextValue "B34" "WEMIXToken"
Set-TextValue "C34" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D34" "2.79"
Set-TextValue "E34" "  -1.18%  "

Set-TextValue "B35" "Monero"
Set-TextValue "C35" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D35" "148.75"
Set-TextValue "E35" "  +3.25%  "

Set-TextValue "D36" "0.0791"
Set-TextValue "E36" "  +3.91%  "

Set-TextValue "D37" "0.116"
Set-TextValue "E37" "  +5.23%  "

Set-TextValue "D38" "0.117"
Set-TextValue "E38" "  +1.87%  "

Set-TextValue "D39" "15.81"
Set-TextValue "E39" "  +3.40%  "

Set-TextValue "D40" "4.09"
Set-TextValue "E40" "  +7.92%  "

Set-TextValue "D41" "0.0307"
Set-TextValue "E41" "  +4.74%  "

Set-TextValue "B42" "NEARProtocol"
Set-TextValue "C42" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D42" "3.36"
Set-TextValue "E42" "  +9.29%  "

Set-TextValue "B43" "Maker"
Set-TextValue "C43" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D43" "2.027.82"
Set-TextValue "E43" "  +13.80%  "

Set-TextValue "B44" "EnergySwap"
Set-TextValue "C44" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D44" "18.89"
Set-TextValue "E44" "  +53.62%  "

Set-TextValue "B45" "FirstDigitalUSD"
Set-TextValue "C45" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D45" "0.999"
Set-TextValue "E45" "  -0.01%  "

Set-TextValue "B46" "BitcoinSV"
Set-TextValue "C46" "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextValue "D46" "92.70"
Set-TextValue "E46" "  +0.49%  "

Set-TextValue "D47" "1.82"
Set-TextValue "E47" "  -2.57%  "

Set-TextValue "B48" "Aave"
Set-TextValue "C48" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D48" "105.10"
Set-TextValue "E48" "  +12.51%  "

Set-TextValue "B49" "FraxShare"
Set-TextValue "C49" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D49" "8.73"
Set-TextValue "E49" "  +12.34%  "

Set-TextValue "B50" "Algorand"
Set-TextValue "C50" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D50" "0.191"
Set-TextValue "E50" "  +4.90%  "

Set-TextValue "B51" "RocketPoolETH"
Set-TextValue "C51" "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue "D51" "2.727.50"
Set-TextValue "E51" "  +10.02%  "

